$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data previously in rows 5-11 (2023-2025 series, without the "subject"
# column) moves up to become the new rows 2-8; the old rows 2-4 (2022 series)
# and the "subject" column are dropped entirely.
# Column K ("1" threshold) and column L ("date") keep their text (string)
# nature - an apostrophe prefix forces Excel to store them as text rather
# than re-interpreting them as a number / date.
$data = @(
    @("AX 11,21", 150, 110, 96,  82, 68, 55, 42, 33, 24, "'16", "'June 2023"),
    @("AY 12,22", 150, 120, 103, 86, 70, 55, 40, 31, 23, "'15", "'June 2023"),
    @("AX 11,21", 150, 116, 100, 84, 68, 53, 38, 30, 23, "'16", "'June 2024"),
    @("AY 12,22", 150, 126, 110, 95, 78, 61, 45, 35, 25, "'15", "'June 2024"),
    @("AX 11,21", 150, 120, 107, 94, 79, 64, 50, 40, 30, "'21", "'June 2025"),
    @("AY 12,22", 150, 129, 114, 100, 84, 69, 54, 45, 37, "'29", "'June 2025"),
    @("AY 12,22", 150, 112, 100, 88, 75, 63, 51, 39, 27, "'15", "'November 2024")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $values[$c]
        if ($col -eq 11 -or $col -eq 12) {
            # Columns K and L were entered with a leading apostrophe so Excel
            # treats them as text instead of a number/date; reset back to
            # the default "Normal" style so no extra quote-prefix styling
            # sticks around on the cell.
            $cell.Style = "Normal"
        }
    }
}

# Drop the now-unused rows 9-11 and the "subject" column (M) entirely so the
# used range shrinks to A1:L8.
$ws.Range("A9:M11").Delete() | Out-Null
$ws.Range("M1:M8").Delete() | Out-Null

$wb.Save()
